# Aggiorno file need_to_buy.xlsx da R
# Refresh the daily forecast data in rows 2-15 (columns A:F) with the
# latest values coming from R (the date window slides forward by one day
# and buy_BEE_MWH / sell_lago_MWH / need_to_buy_MW are recomputed).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2,  45917, 5116.40599940878, 4934.79738876665, 5496, 6531.204028,  35.5664757232445),
    @(3,  45918, 5101.32827217861, 4931.23403399571, 2952, 6512.511375,  141.267380700712),
    @(4,  45919, 5103.30962649499, 4198.19160845483, 2952, 6516.650959,  110.813872539993),
    @(5,  45920, 1207.12767057941, 2150.66626381747, 2952, 2319.790459,  12.9720438432523),
    @(6,  45921, 1096.63641792581, 2126.30773980657, 2952, 2202.28905,   11.6650154950316),
    @(7,  45922, 5230.43751372964, 4813.28820847601, 2952, 6710.218759,  139.211227239432),
    @(8,  45923, 5230.43751372964, 5019.97549315699, 2952, 6710.218759,  147.823197434473),
    @(9,  45924, 5230.43751372964, 5044.81256277744, 2952, 6710.218759,  148.858075335325),
    @(10, 45925, 5230.43751372964, 5044.47008168497, 2952, 6710.218759,  148.843805289806),
    @(11, 45926, 5230.43751372964, 4344.3346962917,  2952, 6710.218759,  119.671497565086),
    @(12, 45927, 1198.04912219854, 2599.21320543046, 2952, 2321.483151,  32.1103014263297),
    @(13, 45928, 1070.29084929054, 2794.54500685562, 2952, 2183.373747,  39.8178293568787),
    @(14, 45929, 5457.44872542722, 5494.98343246843, 2952, 7031.626518,  171.548384376717),
    @(15, 45930, 5457.44872542722, 5770.21147495849, 2952, 7031.626518,  183.016219480469)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Range("A$r").Value = $row[1]
    $ws.Range("B$r").Value = $row[2]
    $ws.Range("C$r").Value = $row[3]
    $ws.Range("D$r").Value = $row[4]
    $ws.Range("E$r").Value = $row[5]
    $ws.Range("F$r").Value = $row[6]
}
